# Foras Promineo CLR FM2 BOM.xlsx - add a second/alternate manufacturer
# ("Manufacturer 2" / "Manufacturer Part Number 2") pair of columns to the
# BOM sheet, with a Murata alternate part populated for the 10u capacitor
# (C5) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foras Promineo CLR PCB")

# New header cells for columns L and M (row 6 is the table header row)
$ws.Range("L6").Value = "Manufacturer 2"
$ws.Range("M6").Value = "Manufacturer Part Number 2"

# New data cells for row 8 (the "C5," / 10u capacitor line item) giving an
# alternate manufacturer / part number
$ws.Range("L8").Value = "Murata Electronics"
$ws.Range("M8").Value = "GRM188C81C106MA73J"

# Size the two new columns like the other "best fit" columns on the sheet
$ws.Columns.Item(12).ColumnWidth = 16.59
$ws.Columns.Item(13).ColumnWidth = 25.75

# Scroll the sheet over a bit and leave the selection where it ended up
# after entering the new data
$ws.Range("C1").Select()
$ws.Range("L18").Select()
